$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 30/12/2025 16:49:20"
$ws1.Range("A3").Value = "Total filas: 430"

$sheet1Rows = @(
    @("", "16:49:09", "16:57", "17_179 Y 38", 8, "LP1912", "30/12/2025"),
    @("", "16:49:09", "16:58", "10_OLMOS", 9, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:01", "16_SANTA ANA", 12, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:05", "11_ETCHEVERRY", 16, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:05", "23_HERNANDEZ", 16, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:11", "10_OLMOS", 22, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:22", "26_HERNANDEZ", 33, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:23", "10_OLMOS", 34, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:25", "84_COLONIA URQUIZA-ESC 49", 36, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:29", "14_ABASTO", 40, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:32", "15_ABASTO", 43, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:35", "23_HERNANDEZ", 46, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:37", "27_EL RETIRO", 48, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:39", "17_ROMERO", 50, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:41", "16_SANTA ANA", 52, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:51", "16_P MOR-167 Y 521", 62, "LP1912", "30/12/2025"),
    @("", "16:49:09", "17:53", "81_EL PELIGRO", 64, "LP1912", "30/12/2025"),
    @("", "16:49:09", "18:05", "17_ROMERO", 76, "LP1912", "30/12/2025"),
    @("", "16:49:09", "18:07", "23_HERNANDEZ", 78, "LP1912", "30/12/2025"),
    @("", "16:49:09", "18:14", "16_SANTA ANA", 85, "LP1912", "30/12/2025"),
    @("", "16:49:09", "18:17", "15_ABASTO", 88, "LP1912", "30/12/2025"),
    @("", "16:49:09", "18:22", "26_HERNANDEZ", 93, "LP1912", "30/12/2025"),
    @("", "16:49:09", "18:25", "14_ABASTO", 96, "LP1912", "30/12/2025")
)

$startRow1 = 409
for ($i = 0; $i -lt $sheet1Rows.Count; $i++) {
    $r = $startRow1 + $i
    $row = $sheet1Rows[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 30/12/2025 16:49:20"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 30/12/2025 16:49:20"
$ws3.Range("A3").Value = "Total filas: 58"

$sheet3Rows = @(
    @("", "30/12/2025", "16:49:19", "16:53", "215B_LP-P MOR-40 Y 115", 4, "L6173"),
    @("", "30/12/2025", "16:49:19", "17:30", "215A_LA PLATA", 41, "L6173"),
    @("", "30/12/2025", "16:49:14", "18:04", "215C_LA PLATA", 75, "L6203")
)

$startRow3 = 57
for ($i = 0; $i -lt $sheet3Rows.Count; $i++) {
    $r = $startRow3 + $i
    $row = $sheet3Rows[$i]
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $ws3.Cells.Item($r, 6).Value = $row[5]
    $ws3.Cells.Item($r, 7).Value = $row[6]
}
